# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right after "2021-Q4" (and before the
#    "总计" summary sheet), populated with the quarter's per-fund holdings.
# 2. Insert a new leading data row into the "总计" summary sheet for the
#    "2022-Q1" quarter, pushing the older quarters down by one row and
#    renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet, positioned between "2021-Q4" and "总计".
#
# NOTE: worksheet handles returned by Worksheets.Item(...) track a *slot
# index*, not a stable sheet identity - once Worksheets.Add() shifts tabs
# around, a handle fetched beforehand can silently start pointing at a
# different sheet. So every handle here is (re-)fetched by name right before
# it's used, and nothing is cached across a structural change (Add/Insert).
# ---------------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4_2021)
$q1.Name = "2022-Q1"

# Match the look & feel of the other per-quarter sheets (bold header row with
# borders, centered column-A index style, page margins, outline flags).
$template = $wb.Worksheets.Item("2021-Q1")   # any per-quarter sheet; used only as a style donor
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A3").PasteSpecial(-4122)

$q1.Outline.SummaryRow    = 1
$q1.Outline.SummaryColumn = 1

$q1.PageSetup.LeftMargin   = 0.75 * 72
$q1.PageSetup.RightMargin  = 0.75 * 72
$q1.PageSetup.TopMargin    = 1    * 72
$q1.PageSetup.BottomMargin = 1    * 72
$q1.PageSetup.HeaderMargin = 0.5  * 72
$q1.PageSetup.FooterMargin = 0.5  * 72

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row 1 - 160416
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'160416"
$q1.Range("C2").Value = "华安标普全球石油指数 (QDII-LOF)"
$q1.Range("D2").Value = "'3.37"
$q1.Range("E2").Value = "'95.08"
$q1.Range("F2").Value = "'9.20"
$q1.Range("G2").Value = "'0.3100"
$q1.Range("H2").Value = 2

# Data row 2 - 000049
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'000049"
$q1.Range("C3").Value = "中银标普全球精选自然资源等权重指数(QDII)"
$q1.Range("D3").Value = "'0.27"
$q1.Range("E3").Value = "'89.72"
$q1.Range("F3").Value = "'1.14"
$q1.Range("G3").Value = "'0.0031"
$q1.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# 2) Insert a new top data row in "总计" for 2022-Q1, shifting the rest down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")   # re-fetch: its slot index shifted from 6 to 7 above
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Restore column-A's bordered/bold index style on the freshly inserted row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.31

# Renumber the index column (A) for the rows that got pushed down.
for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally-active tab (Worksheets.Add activates the new sheet).
$wb.Worksheets.Item(1).Activate()
